# Weekly refresh of the Membrillo / Femacal de La Calera sheet.
#
# The data block (rows 49..82) holds one row per quality ("Especial",
# "Extra (doble especial)", "Primera", sometimes "Segunda") per weekly
# report date, newest week first. A new weekly report is published, so:
#   1. every existing row in the block shifts down by 3 rows
#      (rows 49..82 -> rows 52..85)
#   2. the 3 newest rows (49..51) get this week's figures
#
# We read the existing block into memory first (so the shift doesn't
# clobber data we still need), then write the shifted copy, then the
# brand-new top rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 49
$endRow = 82
$shift = 3
$numCols = 20

# --- 1) snapshot the existing block --------------------------------------
$buffer = @()
for ($r = $startRow; $r -le $endRow; $r++) {
    $rowvals = @()
    for ($c = 1; $c -le $numCols; $c++) {
        $rowvals += $ws.Cells.Item($r, $c).Value2
    }
    $buffer += (,$rowvals)
}

# --- 2) write the snapshot back out, shifted down by 3 rows --------------
for ($i = 0; $i -lt $buffer.Length; $i++) {
    $destRow = $startRow + $shift + $i
    $rowvals = $buffer[$i]
    for ($c = 1; $c -le $numCols; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $rowvals[$c - 1]
    }
}

# --- 3) new rows 83-85 need the date NumberFormat that column D carries --
$ws.Cells.Item(83, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(84, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(85, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- 4) this week's new figures land in rows 49-51 ------------------------
# Common columns for the whole block (market/product identity).
$common = @{
    1  = 3
    2  = "Femacal de La Calera"
    3  = "Coquimbo"
    5  = 5
    6  = "Fruta"
    7  = 100104
    8  = "Frutos de pepita"
    9  = 100104003
    10 = "Membrillo"
    11 = "Champion"
    17 = "`$/caja 18 kilos empedrada"
    18 = "Región de O'Higgins"
    20 = 18
}

function Set-WeekRow {
    param($row, $date, $quality, $volume, $price, $priceKg)

    foreach ($col in $common.Keys) {
        $ws.Cells.Item($row, $col).Value = $common[$col]
    }
    $ws.Cells.Item($row, 4).Value = $date
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 12).Value = $quality
    $ws.Cells.Item($row, 13).Value = $volume
    $ws.Cells.Item($row, 14).Value = $price
    $ws.Cells.Item($row, 15).Value = $price
    $ws.Cells.Item($row, 16).Value = $price
    $ws.Cells.Item($row, 19).Value = $priceKg
}

Set-WeekRow 49 44664 "Especial"                70 14000 778
Set-WeekRow 50 44664 "Extra (doble especial)"   65 16000 889
Set-WeekRow 51 44664 "Primera"                  70 12000 667
